# Reorder the "Recorded By" (column G) comma-separated author lists so that
# a real user email ("dnasr281@gmail.com") is listed first, and failing that
# a "backup@backdoor.com" address is listed first - matching upstream's
# canonical ordering fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $newParts = $null

    if (($parts -contains "dnasr281@gmail.com") -and ($parts[0] -ne "dnasr281@gmail.com")) {
        $rest = $parts | Where-Object { $_ -ne "dnasr281@gmail.com" }
        $newParts = @("dnasr281@gmail.com") + $rest
    }
    elseif (($parts -contains "backup@backdoor.com") -and ($parts[0] -ne "backup@backdoor.com")) {
        $rest = $parts | Where-Object { $_ -ne "backup@backdoor.com" }
        $newParts = @("backup@backdoor.com") + $rest
    }

    if ($null -ne $newParts) {
        $cell.Value = ($newParts -join ", ")
    }
}
